# NPC.xlsx: add a new "Height" property/column for every NPC row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column AB (right after the existing "SkillIDRef" column AA)
$ws.Range("AB1").Value = "Height"

# Give the new column the same width as its neighbour ("SkillIDRef", column AA)
$ws.Columns.Item(28).ColumnWidth = $ws.Columns.Item(27).ColumnWidth

# Every NPC data row (2 through 21) gets a Height value of 2
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 28).Value = 2
}

# Leave the selection on the newly added column, matching the authored workbook
$ws.Range("AB2:AB21").Select() | Out-Null
